# Order_Dummyline_Reconciled.xlsx - reconcile the Sales Order Identifier
# (S2, and the matching "Previous Doc" in AX2) plus the Financial Close
# Calendar Date (T2) against the latest source extract.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# S2 ("Sales Order Identifier") is walked through the candidate order
# numbers pulled from the refreshed source extract before landing on the
# reconciled value.
$ws.Range("S2").Value = "6749007363"
$ws.Range("S2").Value = "7014343960"
$ws.Range("S2").Value = "6751153102"
$ws.Range("S2").Value = "8735182633"
$ws.Range("S2").Value = "4501109159"
$ws.Range("S2").Value = "9842066128"

# AX2 ("Previous Doc") is reconciled to the same order number as S2.
$ws.Range("AX2").Value = "9842066128"

# T2 ("Financial Close Calendar Date DD MM YYYY Code") is stored as plain
# text in this sheet (not an Excel date serial), so stage the new value in
# a scratch cell formatted as text, copy it, and paste-special just the
# value into T2. That preserves T2's original number format/style while
# avoiding Excel auto-converting "01/08/2016" into a date serial. The
# scratch cell (an already-blank cell inside the sheet's used range, so we
# don't disturb the sheet dimensions) is fully cleared afterwards.
$scratch = $ws.Range("AS2")
$scratch.NumberFormat = "@"
$scratch.Value = "01/08/2016"
$scratch.Copy()
$ws.Range("T2").PasteSpecial(-4163)
$scratch.Clear()
